# ---------------------------------------------------------------------------
# feat: support headers containing non-unique keys (#171)
#
# Adds a new "Rebates-Purchases" worksheet (sheet3) with duplicate
# "Rebates"/"Purchases" headers repeated per-month, to exercise the
# xlsx-extractor's handling of non-unique header keys, and nudges a couple
# of cosmetic workbook/sheet view properties.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---- Add the new worksheet as the last tab -------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Rebates-Purchases"

# ---- Row 1: month-end dates, one pair of columns per month ----------------
$ws.Range("C1").Value = 44949
$ws.Range("D1").Value = 44949
$ws.Range("E1").Value = 44980
$ws.Range("F1").Value = 44980
$ws.Range("G1").Value = 45008
$ws.Range("H1").Value = 45008
$ws.Range("I1").Value = 45039
$ws.Range("J1").Value = 45039
$ws.Range("K1").Value = 45069
$ws.Range("L1").Value = 45069
$ws.Range("M1").Value = 45100
$ws.Range("N1").Value = 45100
$ws.Range("O1").Value = 45130
$ws.Range("P1").Value = 45130
$ws.Range("Q1").Value = 45161
$ws.Range("R1").Value = 45161
$ws.Range("S1").Value = 45192
$ws.Range("T1").Value = 45192
$ws.Range("U1").Value = 45222
$ws.Range("V1").Value = 45222
$ws.Range("W1").Value = 45253
$ws.Range("X1").Value = 45253
$ws.Range("Y1").Value = 45283
$ws.Range("Z1").Value = 45283
$ws.Range("C1:Z1").NumberFormat = "d-mmm"

# ---- Row 2: column headers (Rebates/Purchases repeat -> non-unique keys) --
$ws.Range("A2").Value = "Name*"
$ws.Range("B2").Value = "Group*"
$ws.Range("C2").Value = "Rebates"
$ws.Range("D2").Value = "Purchases"
$ws.Range("E2").Value = "Rebates"
$ws.Range("F2").Value = "Purchases"
$ws.Range("G2").Value = "Rebates"
$ws.Range("H2").Value = "Purchases"
$ws.Range("I2").Value = "Rebates"
$ws.Range("J2").Value = "Purchases"
$ws.Range("K2").Value = "Rebates"
$ws.Range("L2").Value = "Purchases"
$ws.Range("M2").Value = "Rebates"
$ws.Range("N2").Value = "Purchases"
$ws.Range("O2").Value = "Rebates"
$ws.Range("P2").Value = "Purchases"
$ws.Range("Q2").Value = "Rebates"
$ws.Range("R2").Value = "Purchases"
$ws.Range("S2").Value = "Rebates"
$ws.Range("T2").Value = "Purchases"
$ws.Range("U2").Value = "Rebates"
$ws.Range("V2").Value = "Purchases"
$ws.Range("W2").Value = "Rebates"
$ws.Range("X2").Value = "Purchases"
$ws.Range("Y2").Value = "Rebates"
$ws.Range("Z2").Value = "Purchases"

# ---- Rows 3-7: sample data -------------------------------------------------
$ws.Range("A3").Value = "John Doe"
$ws.Range("B3").Value = "Group A"
$ws.Range("C3").Value = 100
$ws.Range("D3").Value = 1000
$ws.Range("E3").Value = 110
$ws.Range("F3").Value = 1100
$ws.Range("G3").Value = 120
$ws.Range("H3").Value = 1200
$ws.Range("I3").Value = 130
$ws.Range("J3").Value = 1300
$ws.Range("K3").Value = 140
$ws.Range("L3").Value = 1400
$ws.Range("M3").Value = 150
$ws.Range("N3").Value = 1500
$ws.Range("O3").Value = 160
$ws.Range("P3").Value = 1600
$ws.Range("Q3").Value = 170
$ws.Range("R3").Value = 1700
$ws.Range("S3").Value = 180
$ws.Range("T3").Value = 1800
$ws.Range("U3").Value = 190
$ws.Range("V3").Value = 1900
$ws.Range("W3").Value = 200
$ws.Range("X3").Value = 2000
$ws.Range("Y3").Value = 210
$ws.Range("Z3").Value = 2100

$ws.Range("A4").Value = "Jane Smith"
$ws.Range("B4").Value = "Group B"
$ws.Range("C4").Value = 200
$ws.Range("D4").Value = 2000
$ws.Range("E4").Value = 210
$ws.Range("F4").Value = 2100
$ws.Range("G4").Value = 220
$ws.Range("H4").Value = 2200
$ws.Range("I4").Value = 230
$ws.Range("J4").Value = 2300
$ws.Range("K4").Value = 240
$ws.Range("L4").Value = 2400
$ws.Range("M4").Value = 250
$ws.Range("N4").Value = 2500
$ws.Range("O4").Value = 260
$ws.Range("P4").Value = 2600
$ws.Range("Q4").Value = 270
$ws.Range("R4").Value = 2700
$ws.Range("S4").Value = 280
$ws.Range("T4").Value = 2800
$ws.Range("U4").Value = 290
$ws.Range("V4").Value = 2900
$ws.Range("W4").Value = 300
$ws.Range("X4").Value = 3000
$ws.Range("Y4").Value = 310
$ws.Range("Z4").Value = 3100

$ws.Range("A5").Value = "David Johnson"
$ws.Range("B5").Value = "Group C"
$ws.Range("C5").Value = 300
$ws.Range("D5").Value = 3000
$ws.Range("E5").Value = 310
$ws.Range("F5").Value = 3100
$ws.Range("G5").Value = 320
$ws.Range("H5").Value = 3200
$ws.Range("I5").Value = 330
$ws.Range("J5").Value = 3300
$ws.Range("K5").Value = 340
$ws.Range("L5").Value = 3400
$ws.Range("M5").Value = 350
$ws.Range("N5").Value = 3500
$ws.Range("O5").Value = 360
$ws.Range("P5").Value = 3600
$ws.Range("Q5").Value = 370
$ws.Range("R5").Value = 3700
$ws.Range("S5").Value = 380
$ws.Range("T5").Value = 3800
$ws.Range("U5").Value = 390
$ws.Range("V5").Value = 3900
$ws.Range("W5").Value = 400
$ws.Range("X5").Value = 4000
$ws.Range("Y5").Value = 410
$ws.Range("Z5").Value = 4100

$ws.Range("A6").Value = "Lisa Adams"
$ws.Range("B6").Value = "Group D"
$ws.Range("C6").Value = 400
$ws.Range("D6").Value = 4000
$ws.Range("E6").Value = 410
$ws.Range("F6").Value = 4100
$ws.Range("G6").Value = 420
$ws.Range("H6").Value = 4200
$ws.Range("I6").Value = 430
$ws.Range("J6").Value = 4300
$ws.Range("K6").Value = 440
$ws.Range("L6").Value = 4400
$ws.Range("M6").Value = 450
$ws.Range("N6").Value = 4500
$ws.Range("O6").Value = 460
$ws.Range("P6").Value = 4600
$ws.Range("Q6").Value = 470
$ws.Range("R6").Value = 4700
$ws.Range("S6").Value = 480
$ws.Range("T6").Value = 4800
$ws.Range("U6").Value = 490
$ws.Range("V6").Value = 4900
$ws.Range("W6").Value = 500
$ws.Range("X6").Value = 5000
$ws.Range("Y6").Value = 510
$ws.Range("Z6").Value = 5100

$ws.Range("A7").Value = "Mary Johnson"
$ws.Range("B7").Value = "Group E"
$ws.Range("C7").Value = 500
$ws.Range("D7").Value = 5000
$ws.Range("E7").Value = 510
$ws.Range("F7").Value = 5100
$ws.Range("G7").Value = 520
$ws.Range("H7").Value = 5200
$ws.Range("I7").Value = 530
$ws.Range("J7").Value = 5300
$ws.Range("K7").Value = 540
$ws.Range("L7").Value = 5400
$ws.Range("M7").Value = 550
$ws.Range("N7").Value = 5500
$ws.Range("O7").Value = 560
$ws.Range("P7").Value = 5600
$ws.Range("Q7").Value = 570
$ws.Range("R7").Value = 5700
$ws.Range("S7").Value = 580
$ws.Range("T7").Value = 5800
$ws.Range("U7").Value = 590
$ws.Range("V7").Value = 5900
$ws.Range("W7").Value = 600
$ws.Range("X7").Value = 6000
$ws.Range("Y7").Value = 610
$ws.Range("Z7").Value = 6100

# ---- Column A sizing (best-fit-ish width for the Name* column) -------------
$ws.Columns.Item(1).ColumnWidth = 11.3

# ---- Leave the new sheet's own view state parked on B2, matching how it
#      would have been left after the data was typed in ---------------------
$ws.Activate()
$ws.Range("B2").Select() | Out-Null

# ---- Switch focus back to the Departments tab and nudge its selection -----
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$ws1.Range("D5").Select() | Out-Null
